$d = $word.ActiveDocument

# Step 1: insert one placeholder paragraph per target run, right before the
# document's original first paragraph (anchored on its distinctive opening text).
$anchorRng = $d.Range(0, 0)
$anchorRng.Find.Execute("1. Убрали", $true, $false, $false, $false, $false, $true, 1, $false, "1`r. Что необходимо сделать, чтобы указать `rUE`r, что мы хотим видеть переменные в `rEditor`r’`re`r, специальный тип `rint`r `rв `rUE`r для кроссплатформенности, параметры макроса, как работают`r2. Как создать дочерний `rBlueprint`r-`rкласс`r3. Второй параметр макроса – за что отвечает`r4. Какую функцию имеет каждый `rActor`r (для его имени)`r`r1. Убрали", 2) | Out-Null

# Step 2: set English (en-US) language on the placeholder paragraphs that hold
# English terms, while each still lives alone in its own paragraph (so the
# LanguageID assignment -- which this host applies paragraph-wide -- only
# touches that single run).
$d.Paragraphs(3).Range.LanguageID = "en-US"
$d.Paragraphs(5).Range.LanguageID = "en-US"
$d.Paragraphs(7).Range.LanguageID = "en-US"
$d.Paragraphs(9).Range.LanguageID = "en-US"
$d.Paragraphs(12).Range.LanguageID = "en-US"
$d.Paragraphs(15).Range.LanguageID = "en-US"
$d.Paragraphs(20).Range.LanguageID = "en-US"

# Step 3: stitch the placeholder paragraphs for each target paragraph back
# together by deleting the paragraph mark between consecutive runs that belong
# to the same target paragraph (run formatting/boundaries are preserved).

$cursor = 1
# target paragraph 1: 13 run(s)
for ($i = 0; $i -lt 12; $i++) {
  $endp = $d.Paragraphs($cursor).Range.End
  $d.Range($endp - 1, $endp).Delete()
}
$cursor = $cursor + 1
# target paragraph 2: 4 run(s)
for ($i = 0; $i -lt 3; $i++) {
  $endp = $d.Paragraphs($cursor).Range.End
  $d.Range($endp - 1, $endp).Delete()
}
$cursor = $cursor + 1
# target paragraph 3: 1 run(s)
$cursor = $cursor + 1
# target paragraph 4: 3 run(s)
for ($i = 0; $i -lt 2; $i++) {
  $endp = $d.Paragraphs($cursor).Range.End
  $d.Range($endp - 1, $endp).Delete()
}
$cursor = $cursor + 1
# target paragraph 5: 1 run(s)
$cursor = $cursor + 1

Write-Host "done"